# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.353.25"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.711.71"
$ws.Range("E3").Value = "  +9.08%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'609.13"
$ws.Range("E5").Value = "  +3.85%  "
$ws.Range("D6").Value = "'174.95"
$ws.Range("E6").Value = "  -4.51%  "
$ws.Range("D7").Value = "3.714.04"
$ws.Range("E7").Value = "  +9.31%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.537"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "'40.43"
$ws.Range("E13").Value = "  +5.33%  "
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "4.326.10"
$ws.Range("E15").Value = "  +9.24%  "
$ws.Range("D16").Value = "3.713.77"
$ws.Range("E16").Value = "  +9.22%  "
$ws.Range("D17").Value = "69.438.36"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'7.54"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'16.65"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'510.90"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "'9.37"
$ws.Range("E22").Value = "  +12.08%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'87.43"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("D26").Value = "'13.30"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +22.34%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  -4.63%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.82"
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("D33").Value = "'31.01"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "'2.16"
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").Value = "'51.16"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").Value = "'43.96"
$ws.Range("E42").Value = "  -10.16%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.74"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "3.071.75"
$ws.Range("E44").Value = "  +4.64%  "
$ws.Range("D45").Value = "'416.89"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").Value = "'2.68"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("D47").Value = "'0.0361"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").Value = "'27.61"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'134.70"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.48"
$ws.Range("E51").Value = "  +1.73%  "
